# Refresh the cryptocurrency price/volume snapshot (Price column D,
# Volume(1h) column E) for rows 2-51 to match the latest scraped values.
# Column D values are prefixed with a leading apostrophe so Excel stores
# them as text (preserving formatting such as trailing zeros and the
# "thousands.dot" style numbers) instead of auto-converting them to
# numeric cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'41.314.87"
$ws.Range("E2").Value = "  -3.74%  "

$ws.Range("D3").Value = "'2.467.58"
$ws.Range("E3").Value = "  -2.72%  "

$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").Value = "'311.62"
$ws.Range("E5").Value = "  -0.25%  "

$ws.Range("D6").Value = "'94.70"
$ws.Range("E6").Value = "  -5.83%  "

$ws.Range("D7").Value = "'0.549"
$ws.Range("E7").Value = "  -3.24%  "

$ws.Range("E8").Value = "  -0.14%  "

$ws.Range("E9").Value = "  -4.62%  "

$ws.Range("D10").Value = "'33.57"
$ws.Range("E10").Value = "  -6.43%  "

$ws.Range("D11").Value = "'0.0782"
$ws.Range("E11").Value = "  -3.13%  "

$ws.Range("E12").Value = "  -1.13%  "

$ws.Range("E13").Value = "  -4.55%  "

$ws.Range("D14").Value = "'2.847.26"
$ws.Range("E14").Value = "  -2.64%  "

$ws.Range("D15").Value = "'2.513.77"
$ws.Range("E15").Value = "  -0.35%  "

$ws.Range("D16").Value = "'14.95"
$ws.Range("E16").Value = "  -3.41%  "

$ws.Range("E17").Value = "  -3.95%  "

$ws.Range("D18").Value = "'41.319.71"
$ws.Range("E18").Value = "  -3.65%  "

$ws.Range("E19").Value = "  -5.56%  "

$ws.Range("D20").Value = "'0.0₃0925"

$ws.Range("E21").Value = "  -9.92%  "

$ws.Range("D22").Value = "'68.57"
$ws.Range("E22").Value = "  -1.85%  "

$ws.Range("D23").Value = "'237.16"
$ws.Range("E23").Value = "  -2.89%  "

$ws.Range("E24").Value = "  -4.62%  "

$ws.Range("E25").Value = "  -0.03%  "

$ws.Range("E26").Value = "  -6.97%  "

$ws.Range("D27").Value = "'24.12"
$ws.Range("E27").Value = "  -6.12%  "

$ws.Range("D28").Value = "'2.22"
$ws.Range("E28").Value = "  -5.14%  "

$ws.Range("D29").Value = "'9.65"
$ws.Range("E29").Value = "  -5.94%  "

$ws.Range("D30").Value = "'36.64"
$ws.Range("E30").Value = "  -6.12%  "

$ws.Range("D31").Value = "'151.89"

$ws.Range("D32").Value = "'5.50"
$ws.Range("E32").Value = "  -6.04%  "

$ws.Range("D33").Value = "'2.66"
$ws.Range("E33").Value = "  -4.51%  "

$ws.Range("D34").Value = "'2.58"
$ws.Range("E34").Value = "  -3.36%  "

$ws.Range("E35").Value = "  -5.61%  "

$ws.Range("E36").Value = "  -3.50%  "

$ws.Range("D37").Value = "'17.14"
$ws.Range("E37").Value = "  -6.85%  "

$ws.Range("E38").Value = "  -5.33%  "

$ws.Range("E39").Value = "  -3.16%  "

$ws.Range("D40").Value = "'0.103"
$ws.Range("E40").Value = "  -7.82%  "

$ws.Range("D41").Value = "'4.24"
$ws.Range("E41").Value = "  +0.90%  "

$ws.Range("E42").Value = "  -0.05%  "

$ws.Range("D43").Value = "'19.74"
$ws.Range("E43").Value = "  -10.53%  "

$ws.Range("D44").Value = "'1.987.48"
$ws.Range("E44").Value = "  -0.81%  "

$ws.Range("E45").Value = "  -4.91%  "

$ws.Range("D46").Value = "'3.03"
$ws.Range("E46").Value = "  -9.53%  "

$ws.Range("D47").Value = "'8.72"
$ws.Range("E47").Value = "  -5.95%  "

$ws.Range("D48").Value = "'2.712.49"
$ws.Range("E48").Value = "  -2.26%  "

$ws.Range("D49").Value = "'69.73"
$ws.Range("E49").Value = "  -4.09%  "

$ws.Range("D50").Value = "'96.50"
$ws.Range("E50").Value = "  -5.07%  "

$ws.Range("E51").Value = "  -7.56%  "
